# Apply the "mA -> A" conversion to column B of the LCR data sheet.
# - Rename header B1 from "I [mA]" to "I [A]"
# - Convert each data value in B2:B20 from milliamps to amps (divide by 1000)
#   and round to 2 decimal places, matching the "rounded" output file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header.
$ws.Range("B1").Value = "I [A]"

# Rows 2-20 hold the measured current values (in mA) that need converting to A.
for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $mA = $cell.Value2
    $A = [Math]::Round($mA / 1000, 2)
    $cell.Value = $A
}
